# Auto-generated Excel COM-interop edit script
# Applies crypto price/volume updates per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.286.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.680.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5274"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.009"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2711"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06477"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07518"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.680.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.537"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5815"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008504"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.324.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.928"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.203"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.009"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.827"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1243"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06564"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.357"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.331"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.604"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.596"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.661"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.035"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6250"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.404"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.751"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.465"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.115.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01625"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8791"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.015"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.831.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.154"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.003"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05280"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.104"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4297"
$ws.Range("D51").Style = "Normal"
